# #fix add event bug #edit footer
#
# 1) Add a new permission row "CanAddEvent" (order 10) at A9:B9, and
#    renumber the existing "order" column (A) back to a clean 1..9
#    sequence (the sheet is kept sorted by column A).
# 2) Add a sample "User" record in D1:K1 (Id, Name, PasswordHash,
#    Country, Email, some numeric code, and two boolean flags).
# 3) Re-point the selection, widen/add a few column widths, and touch
#    up the page setup (footer / print orientation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new permission row -----------------------------------------------
$ws.Range("A9").Value = 10
$ws.Range("B9").Value = "CanAddEvent"

# Re-sort A1:A9 (carries the paired B column with it) so the "order"
# column goes back to a clean ascending sequence.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1"))
$ws.Sort.SetRange($ws.Range("A1:A9"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# --- sample user record --------------------------------------------
$ws.Range("D1").Value = 1
$ws.Range("E1").Value = "Ahmed"
$ws.Range("F1").Value = "202cb962ac59075b964b07152d234b70"
$ws.Range("G1").Value = "egypt"
$ws.Range("H1").Value = "a@mail.com"
$ws.Range("I1").Value = 1148177915
$ws.Range("J1").Value = $true
$ws.Range("K1").Value = $false

# match the existing "vertical center" cell style (style index 1) used
# by the rest of the sheet
$ws.Range("D1:K1").VerticalAlignment = -4108

# --- column widths ------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 2.166666666666667
$ws.Columns.Item(6).ColumnWidth = 32.498697916666664
$ws.Columns.Item(8).ColumnWidth = 10.498697916666666
$ws.Columns.Item(9).ColumnWidth = 10.166666666666666

# --- selection / view ----------------------------------------------
$null = $ws.Range("I2").Select()

# --- page setup (footer / print orientation) ------------------------
$ws.PageSetup.Orientation = 1
